$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (strikeouts) values per row, replacing the previous Strike# values.
$kValues = @{
    2 = 4
    3 = 1
    5 = 3
    6 = 2
    7 = 1
    8 = 1
    9 = 1
    10 = 0
    11 = 3
    12 = 1
    13 = 3
    14 = 1
    15 = 1
    16 = 2
    17 = 3
    18 = 2
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 2
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 1
    32 = 3
    33 = 0
    34 = 2
    35 = 1
    36 = 1
    37 = 1
    38 = 2
    39 = 1
    40 = 2
    41 = 1
    42 = 0
    43 = 2
    44 = 0
    45 = 1
    46 = 1
    47 = 2
    48 = 1
    49 = 2
    50 = 2
    51 = 3
    52 = 1
    53 = 0
    54 = 1
    55 = 2
    56 = 1
    57 = 2
    58 = 0
    59 = 1
    60 = 1
    61 = 0
    62 = 0
    63 = 1
    64 = 0
    65 = 1
    66 = 1
    67 = 0
    68 = 3
    69 = 2
    70 = 0
    71 = 1
    72 = 1
    73 = 1
    74 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}
